# The "Ink / Varnish" table lists one row per ink color (Yellow, Black,
# Cyan, Magenta) in columns D (color/process description) and H (material).
# The edit reorders these rows so that Cyan comes first (row 2), pushing
# Yellow down to row 3 and Black down to row 4; Magenta (row 5) is unchanged.
#
#   Before: Row2=Yellow, Row3=Black, Row4=Cyan, Row5=Magenta
#   After:  Row2=Cyan,   Row3=Yellow, Row4=Black, Row5=Magenta

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for the two columns that vary
# per ink row.
$d2 = $ws.Range("D2").Value2
$h2 = $ws.Range("H2").Value2
$d3 = $ws.Range("D3").Value2
$h3 = $ws.Range("H3").Value2
$d4 = $ws.Range("D4").Value2
$h4 = $ws.Range("H4").Value2

# Row 2 becomes what Row 4 (Cyan) used to be.
$ws.Range("D2").Value = $d4
$ws.Range("H2").Value = $h4

# Row 3 becomes what Row 2 (Yellow) used to be.
$ws.Range("D3").Value = $d2
$ws.Range("H3").Value = $h2

# Row 4 becomes what Row 3 (Black) used to be.
$ws.Range("D4").Value = $d3
$ws.Range("H4").Value = $h3

# Row 5 (Magenta) is left untouched.
